# master-reg_center_user_machine_h.xlsx: append 9 more test rows, move the
# selection, and switch the print page to portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already holds rows 2-21 (regcntr_id 10001..10015, machine_id
# 10001..10020). Continue the same sequence for 9 more rows (22-30):
# regcntr_id restarts at 10002, usr_id/machine_id keep climbing.
$rows = @(
    @(10002, 110021, 10021),
    @(10003, 110022, 10022),
    @(10004, 110023, 10023),
    @(10005, 110024, 10024),
    @(10006, 110025, 10025),
    @(10007, 110026, 10026),
    @(10008, 110027, 10027),
    @(10009, 110028, 10028),
    @(10010, 110029, 10029)
)

$r = 22
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]   # regcntr_id
    $ws.Cells.Item($r, 2).Value = $row[1]   # usr_id
    $ws.Cells.Item($r, 3).Value = $row[2]   # machine_id
    $ws.Cells.Item($r, 4).Value = "eng"         # lang_code
    $ws.Cells.Item($r, 5).Value = $true         # is_active
    $ws.Cells.Item($r, 6).Value = "superadmin"  # cr_by
    $ws.Cells.Item($r, 7).Value = "now()"       # cr_dtimes
    $ws.Cells.Item($r, 8).Value = "now()"       # eff_dtimes
    $r++
}

# The saved workbook's last active cell was F14.
$ws.Range("F14").Select()

# Switch the sheet's print setup to portrait (xlPortrait = 1).
$ws.PageSetup.Orientation = 1
